# Generate Report for handback
# Refresh the "Correspond Handoff Datetime" (D2) and "Correspond Handback
# DateTime" (G2) values for the first file row on each language sheet, to
# reflect a newly generated handback report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-01-11 03:00:27"
$zhcn.Range("G2").Value = "2016-01-11 03:01:36"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-01-11 03:00:43"
$dede.Range("G2").Value = "2016-01-11 03:02:00"
